$wb = $excel.ActiveWorkbook

# --- Vendors sheet: reduce remaining quantities for two rows ---
$wsVendors = $wb.Worksheets.Item("Vendors")
$wsVendors.Range("D6").Value = 95
$wsVendors.Range("D8").Value = 95

# --- Clients sheet: append two new client rows ---
$wsClients = $wb.Worksheets.Item("Clients")
$wsClients.Range("A10").Value = "Imara Gimura"
$wsClients.Range("B10").Value = "imarasnakeeyes@EBI.com"
$wsClients.Range("A11").Value = "Vincent Pryor"
$wsClients.Range("B11").Value = "DrPryor@EBI.com"

# --- Expenses sheet: append two new expense rows ---
$wsExpenses = $wb.Worksheets.Item("Expenses")
$wsExpenses.Range("A6").Value = "Imara Gimura"
$wsExpenses.Range("B6").Value = "Zelda BOTW (5% Discount)"
$wsExpenses.Range("C6").Value = 5
$wsExpenses.Range("D6").Value = 189.95249999999999

$wsExpenses.Range("A7").Value = "Vincent Pryor"
$wsExpenses.Range("B7").Value = "Samsung S9000"
$wsExpenses.Range("C7").Value = 5
$wsExpenses.Range("D7").Value = 4499.95
